$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.834.66"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.733.09"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5143"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "1.749.35"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07027"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6428"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "25.827.45"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006604"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "1.972.71"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.143"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.687"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.125"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("E27").Value = "  +3.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.794"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08318"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.681"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.418"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9804"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01579"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3830"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.962"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05392"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1120"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.641"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.77%  "
